# Append the new daily GSC export row (2025-11-15) to the "Chart" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$row = 42

# Column A holds the date as plain text (e.g. "2025-11-14"), not a real
# Excel date value, so we temporarily force a text number format to stop
# Excel from auto-converting the string into a date serial number, then
# clear the formatting back to the sheet's default ("General") style so
# the new cell matches the style of the existing rows.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-11-15"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 35
